# Big results push, waiting on protonet to finish
#
# Update the meanrank table:
#   - row 2 relabeled to the new "GNN-MT-O" variant, with a refreshed value
#   - RF / GNN-MT rows keep their position but get refreshed values
#   - two brand-new rows are appended for the "-O" (ordinal) variants of
#     GNN-MT and PN, and the PN row itself gets a refreshed value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: now the GNN-MT-O entry
$ws.Range("A2").Value = "8_train (GNN-MT-O) val delta-auprc"
$ws.Range("B2").Value = 3.633333333333333

# Row 3: RF entry, refreshed value
$ws.Range("A3").Value = "8_train (RF) val delta-auprc"
$ws.Range("B3").Value = 3.274193548387097

# Row 4: GNN-MT entry, refreshed value
$ws.Range("A4").Value = "8_train (GNN-MT) val delta-auprc"
$ws.Range("B4").Value = 3.066666666666667

# Row 5 (new): PN-O entry - copy formatting from row 4's label cell first
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "8_train (PN-O) val delta-auprc"
$ws.Range("B5").Value = 2.461290322580645

# Row 6 (new): PN entry - copy formatting from row 4's label cell first
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "8_train (PN) val delta-auprc"
$ws.Range("B6").Value = 2.267741935483871

$excel.CutCopyMode = $false
